# Fixed #476 Moving from Apache POI 4.1.0 to 5.2.3.
#
# The generated table in this document had every "label" run (the bold
# black 12pt runs used for the Sirius representation tree labels) written
# with a run-properties block in the order b/i/strike/color/sz and with
# literal "true"/"false" boolean spellings. After the POI upgrade the same
# semantic formatting (Bold=True, Italic=False, StrikeThrough=False,
# Color=black, Size=12) is re-emitted (order sz/b/i/strike/color, on/off
# spellings). We reapply the same formatting via the Word object model so
# every one of those runs keeps its Bold/Italic/StrikeThrough/Color/Size
# values.

$d = $word.ActiveDocument

$total = $d.Content.End

$i = 0
while ($i -lt $total) {
    $probe = $d.Range($i, $i + 1)
    if ($probe.Font.Bold -eq -1) {
        $segStart = $i
        $j = $i
        while ($j -lt $total) {
            $probe2 = $d.Range($j, $j + 1)
            if ($probe2.Font.Bold -ne -1) {
                break
            }
            $j = $j + 1
        }
        $segEnd = $j

        $rng = $d.Range($segStart, $segEnd)
        $rng.Font.Bold = -1
        $rng.Font.Italic = 0
        $rng.Font.StrikeThrough = 0
        $rng.Font.Color = 0
        $rng.Font.Size = 12

        $i = $segEnd
    } else {
        $i = $i + 1
    }
}
